# Correct typos in the "Comments" column (B) of the PyTorch support matrix:
# add missing trailing periods, and fix row 27 (Squeeze) height which was
# oversized for its wrapped two-line text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Conv2d/MaxPool2d dilation note
$ws.Range("B6").Value = "Only dilation=1."

# Row 7 - AvgPool2d dilation/count_include_pad note
$ws.Range("B7").Value = "Only dilation=1, count_include_pad=1."

# Row 9 - interpolate note
$ws.Range("B9").Value = "See Supported operations (ONNX): resize."

# Row 10 - Upsample note (two lines)
$ws.Range("B10").Value = "Only mode ""nearest"" and ""linear"".`nOnly scales=[2,2]."

# Row 11 - add alpha note
$ws.Range("B11").Value = "Only alpha=1."

# Row 12 - sub alpha note (duplicate text of row 11, now de-duplicated by Excel)
$ws.Range("B12").Value = "Only alpha=1."

# Row 13 - mul note
$ws.Range("B13").Value = "Only constant multiplication."

# Row 14 - div note
$ws.Range("B14").Value = "Only constant division."

# Row 15 - cat note
$ws.Range("B15").Value = "Only along channel axis."

# Row 21 - Softplus beta note
$ws.Range("B21").Value = "Only beta=1."

# Row 25 - Reshape note
$ws.Range("B25").Value = "Only channel-wise flatten and before fully connected layer or Conv w/ 1x1 kernel."

# Row 26 - Transpose note
$ws.Range("B26").Value = "Only before fully connected layer."

# Row 27 - Squeeze note (two lines) + row height correction (was oversized at 66)
$ws.Range("B27").Value = "Only when resulting tensor has 2D shape.`nSqueeze along batch axis is unsupported."
$ws.Rows.Item(27).RowHeight = 33

# Row 28 - Flatten note (same text as row 25)
$ws.Range("B28").Value = "Only channel-wise flatten and before fully connected layer or Conv w/ 1x1 kernel."

# Update the saved cursor/selection position to reflect the author's final
# cell selection after editing.
$ws.Range("B33").Select()
